$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values per the new run log data.
$ws.Range("C2:C6").Value = 8552
$ws.Range("C7:C17").Value = 8441
$ws.Range("C18:C21").Value = 8225
$ws.Range("C22:C24").Value = 8042
$ws.Range("C25:C29").Value = 7947
$ws.Range("C30:C252").Value = 7310
